# Generate Report for Handback
# Refresh the handback-status report with a new source/handoff pair:
#   59199e9b-216c-4568-9c85-c61bf9ca802a  ->  5404d842-ec33-4928-b046-304a8d7de8e3
#   cc515795-00a0-4335-9ec6-ff4d30aa0e5e  ->  ffffb5aa4b5f-4ef0-4a03-a969-b94bc9e7e4ce
# and new handoff/handback hashes + timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("G2").Value = "2016-08-20 23:06:20"
$ov.Range("G3").Value = "2016-08-20 23:06:20"

# Hyperlinks (also sets the displayed cell text for B2/B3).
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/039a68254f545540d1af906768275cb8c1b499ae/e2e/5404d842-ec33-4928-b046-304a8d7de8e3.md", [type]::Missing, [type]::Missing, "e2e\5404d842-ec33-4928-b046-304a8d7de8e3.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/039a68254f545540d1af906768275cb8c1b499ae/e2e/ffffb5aa4b5f-4ef0-4a03-a969-b94bc9e7e4ce.md", [type]::Missing, [type]::Missing, "e2e\ffffb5aa4b5f-4ef0-4a03-a969-b94bc9e7e4ce.md") | Out-Null

# Plain text cells (not hyperlinked) in A/B still need the raw file name refreshed.
$ov.Range("A2").Value = "5404d842-ec33-4928-b046-304a8d7de8e3.md"
$ov.Range("A3").Value = "ffffb5aa4b5f-4ef0-4a03-a969-b94bc9e7e4ce.md"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("G2").Value = "5404d842-ec33-4928-b046-304a8d7de8e3.8c79448fc3292a36973aa945304b205c14e8dac8.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-20 23:06:15"
$zh.Range("J2").Value = "5404d842-ec33-4928-b046-304a8d7de8e3.8c79448fc3292a36973aa945304b205c14e8dac8.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-20 23:06:32"

$zh.Range("G3").Value = "5404d842-ec33-4928-b046-304a8d7de8e3.8c79448fc3292a36973aa945304b205c14e8dac8.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-20 23:06:15"
$zh.Range("J3").Value = "5404d842-ec33-4928-b046-304a8d7de8e3.8c79448fc3292a36973aa945304b205c14e8dac8.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-20 23:06:32"

# Hyperlinks (also sets the displayed cell text for A2/I2/A3/I3).
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/039a68254f545540d1af906768275cb8c1b499ae/e2e/59199e9b-216c-4568-9c85-c61bf9ca802a.md", [type]::Missing, [type]::Missing, "5404d842-ec33-4928-b046-304a8d7de8e3.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f4e7f745e1c079eadb148503bcf0940b5ab70c5e/e2e/59199e9b-216c-4568-9c85-c61bf9ca802a.md", [type]::Missing, [type]::Missing, "5404d842-ec33-4928-b046-304a8d7de8e3.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/039a68254f545540d1af906768275cb8c1b499ae/e2e/cc515795-00a0-4335-9ec6-ff4d30aa0e5e.md", [type]::Missing, [type]::Missing, "ffffb5aa4b5f-4ef0-4a03-a969-b94bc9e7e4ce.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f4e7f745e1c079eadb148503bcf0940b5ab70c5e/e2e/cc515795-00a0-4335-9ec6-ff4d30aa0e5e.md", [type]::Missing, [type]::Missing, "ffffb5aa4b5f-4ef0-4a03-a969-b94bc9e7e4ce.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("G2").Value = "5404d842-ec33-4928-b046-304a8d7de8e3.8c79448fc3292a36973aa945304b205c14e8dac8.de-de.xlf"
$de.Range("H2").Value = "2016-08-20 23:06:20"
$de.Range("J2").Value = "5404d842-ec33-4928-b046-304a8d7de8e3.8c79448fc3292a36973aa945304b205c14e8dac8.de-de.xlf"
$de.Range("K2").Value = "2016-08-20 23:06:38"

$de.Range("G3").Value = "5404d842-ec33-4928-b046-304a8d7de8e3.8c79448fc3292a36973aa945304b205c14e8dac8.de-de.xlf"
$de.Range("H3").Value = "2016-08-20 23:06:20"
$de.Range("J3").Value = "5404d842-ec33-4928-b046-304a8d7de8e3.8c79448fc3292a36973aa945304b205c14e8dac8.de-de.xlf"
$de.Range("K3").Value = "2016-08-20 23:06:38"

# Hyperlinks (also sets the displayed cell text for A2/I2/A3/I3).
$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/039a68254f545540d1af906768275cb8c1b499ae/e2e/59199e9b-216c-4568-9c85-c61bf9ca802a.md", [type]::Missing, [type]::Missing, "5404d842-ec33-4928-b046-304a8d7de8e3.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/87964a9011d7fcd7a92ba261d527d2519c902cb6/e2e/59199e9b-216c-4568-9c85-c61bf9ca802a.md", [type]::Missing, [type]::Missing, "5404d842-ec33-4928-b046-304a8d7de8e3.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/039a68254f545540d1af906768275cb8c1b499ae/e2e/cc515795-00a0-4335-9ec6-ff4d30aa0e5e.md", [type]::Missing, [type]::Missing, "ffffb5aa4b5f-4ef0-4a03-a969-b94bc9e7e4ce.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/87964a9011d7fcd7a92ba261d527d2519c902cb6/e2e/cc515795-00a0-4335-9ec6-ff4d30aa0e5e.md", [type]::Missing, [type]::Missing, "ffffb5aa4b5f-4ef0-4a03-a969-b94bc9e7e4ce.md") | Out-Null
